{"js": "// Update the date line and the 25 three-digit-by-one-digit multiplication\n// problems in the practice sheet to the next day's worksheet values.\nconst replacements = [\n  [\"2025-06-19 Thursday\", \"2025-06-20 Friday\"],\n  [\"460\u00d77=\", \"681\u00d79=\"],\n  [\"746\u00d74=\", \"234\u00d73=\"],\n  [\"743\u00d73=\", \"306\u00d74=\"],\n  [\"867\u00d73=\", \"868\u00d74=\"],\n  [\"753\u00d79=\", \"474\u00d77=\"],\n  [\"746\u00d77=\", \"611\u00d72=\"],\n  [\"909\u00d74=\", \"764\u00d72=\"],\n  [\"845\u00d75=\", \"403\u00d77=\"],\n  [\"681\u00d75=\", \"523\u00d77=\"],\n  [\"614\u00d72=\", \"979\u00d78=\"],\n  [\"649\u00d73=\", \"648\u00d76=\"],\n  [\"389\u00d77=\", \"754\u00d77=\"],\n  [\"265\u00d75=\", \"117\u00d76=\"],\n  [\"673\u00d72=\", \"900\u00d78=\"],\n  [\"664\u00d75=\", \"135\u00d76=\"],\n  [\"563\u00d79=\", \"193\u00d76=\"],\n  [\"916\u00d76=\", \"671\u00d73=\"],\n  [\"374\u00d78=\", \"400\u00d78=\"],\n  [\"479\u00d74=\", \"572\u00d75=\"],\n  [\"973\u00d77=\", \"995\u00d73=\"],\n  [\"495\u00d77=\", \"298\u00d72=\"],\n  [\"411\u00d78=\", \"600\u00d73=\"],\n  [\"811\u00d76=\", \"943\u00d77=\"],\n  [\"159\u00d78=\", \"612\u00d75=\"],\n  [\"851\u00d77=\", \"892\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 three-digit-by-one-digit multiplication\n# problems in the practice sheet to the next day's worksheet values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-19 Thursday\", \"2025-06-20 Friday\"),\n    @(\"460\u00d77=\", \"681\u00d79=\"),\n    @(\"746\u00d74=\", \"234\u00d73=\"),\n    @(\"743\u00d73=\", \"306\u00d74=\"),\n    @(\"867\u00d73=\", \"868\u00d74=\"),\n    @(\"753\u00d79=\", \"474\u00d77=\"),\n    @(\"746\u00d77=\", \"611\u00d72=\"),\n    @(\"909\u00d74=\", \"764\u00d72=\"),\n    @(\"845\u00d75=\", \"403\u00d77=\"),\n    @(\"681\u00d75=\", \"523\u00d77=\"),\n    @(\"614\u00d72=\", \"979\u00d78=\"),\n    @(\"649\u00d73=\", \"648\u00d76=\"),\n    @(\"389\u00d77=\", \"754\u00d77=\"),\n    @(\"265\u00d75=\", \"117\u00d76=\"),\n    @(\"673\u00d72=\", \"900\u00d78=\"),\n    @(\"664\u00d75=\", \"135\u00d76=\"),\n    @(\"563\u00d79=\", \"193\u00d76=\"),\n    @(\"916\u00d76=\", \"671\u00d73=\"),\n    @(\"374\u00d78=\", \"400\u00d78=\"),\n    @(\"479\u00d74=\", \"572\u00d75=\"),\n    @(\"973\u00d77=\", \"995\u00d73=\"),\n    @(\"495\u00d77=\", \"298\u00d72=\"),\n    @(\"411\u00d78=\", \"600\u00d73=\"),\n    @(\"811\u00d76=\", \"943\u00d77=\"),\n    @(\"159\u00d78=\", \"612\u00d75=\"),\n    @(\"851\u00d77=\", \"892\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
